# "hours update source add"
# Add three new timesheet rows (40-42) to Sheet1 with dates, hours and
# comments, and refresh the sheet's view/selection + column widths to
# match the post-edit layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New data rows -------------------------------------------------------
# Row 40: 2010-03-23, 2.5 hrs, Group Meeting
$ws.Range("A40").Value = 40260
$ws.Range("B40").Value = 2.5
$ws.Range("C40").Value = "Group Meeting"

# Row 41: 2010-03-23, 1 hr, Weekly Meeting
$ws.Range("A41").Value = 40260
$ws.Range("B41").Value = 1
$ws.Range("C41").Value = "Weekly Meeting"

# Row 42: 2010-03-26, 0.5 hrs, Skype Meeting
$ws.Range("A42").Value = 40263
$ws.Range("B42").Value = 0.5
$ws.Range("C42").Value = "Skype Meeting"

# --- View / layout bookkeeping -------------------------------------------
# Move the active selection down to the next empty row (A43)
[void]$ws.Range("A43").Select()

# Column widths shrink very slightly in the source edit
$ws.Columns("A").ColumnWidth = 10.75
$ws.Columns("B").ColumnWidth = 6.45
